$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "2019年5月9日23:05:28"
$ws.Range("B44").Value = "周四"
$ws.Range("C44").Value = "dao service bug修正，完善"
$ws.Range("D44").Value = "10:30--12:00 & 13:30--15:00"

$ws.Range("C45").Value = "JavaEE上机完善"
$ws.Range("D45").Value = "10:30--12:00 & 13:30--15:00"

$ws.Range("C46").Value = "交互控制 Controller"
$ws.Range("D46").Value = "16:30--19:00"

$ws.Range("D46").Select()
